# Updated results to April 22, removed France, changed file names of output
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that disappear entirely (old USA/NYC=row7 data stays,
# but old row 8 "Spain" and row 9 "Italy" shift up as rows 6/7 already hold
# new data below, so we remove the two now-unused trailing rows from the
# bottom after rewriting the top). Also "France" (old row 6) is dropped.
# Simplest approach: delete old row 6 (France) and old row 2 (USA/All) since
# the surviving rows shift into their place, then overwrite remaining cells
# with the refreshed values from the new data pull (and clear the retired
# H/I helper columns on each data row).

$ws.Rows.Item(6).Delete()  # France
$ws.Rows.Item(2).Delete()  # USA / All (old, superseded by refreshed pull below)

# After the two deletions the sheet now has exactly the six data rows
# (rows 2-7) required by the new dataset. Overwrite them with the refreshed
# figures and drop the now-unused relAgeDE / relRateDE columns (H:I).

# Column C keeps its existing date-format style (s="2" -> yyyy-mm-dd); the
# underlying cell value is the usual 1900-epoch date serial number, so we
# just assign the serial directly (2020-04-22 -> 43943, 2020-02-11 -> 43872).
$data = @(
    @{ Row = 2; A = "SouthKorea"; B = "All"; C = 43943; D = 0.02226482139517487;  E = -0.02226482139517487 },
    @{ Row = 3; A = "China";      B = "All"; C = 43872; D = 0.02290025071633238;  E = -0.02290025071633238 },
    @{ Row = 4; A = "Germany";    B = "All"; C = 43943; D = 0.03807511017180802;  E = -0.03807511017180802 },
    @{ Row = 5; A = "USA";        B = "All"; C = 43943; D = 0.05598346774242222;  E = -0.05598346774242222 },
    @{ Row = 6; A = "Spain";      B = "All"; C = 43943; D = 0.1040122728013148;   E = -0.1040122728013148 },
    @{ Row = 7; A = "Italy";      B = "All"; C = 43943; D = 0.1300063316640764;   E = -0.1300063316640764 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).ClearContents()
    $ws.Cells.Item($r, 9).ClearContents()
}
